# Delete the "Rectangle: Rounded Corners 29" shape (id 28) and the
# "Connector: Elbow 31" connector (id 29) from the slide with SlideID 386
# (the "%cr3" callout box and its elbow connector into Rectangle 21).

$p = $ppt.ActivePresentation

# Locate the slide by its persistent SlideID (386) rather than a fixed
# positional index, so the script is robust to slide-order differences.
$targetSlide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $candidate = $p.Slides.Item($si)
    if ($candidate.SlideID -eq 386) {
        $targetSlide = $candidate
        break
    }
}

$idsToDelete = @(28, 29)

# Walk the shape collection back-to-front and delete by the shapes'
# stable .Id property (not by collection index, which shifts as shapes
# are removed, and not by .Name, since this slide has two connectors
# both named "Connector: Elbow 31").
for ($i = $targetSlide.Shapes.Count; $i -ge 1; $i--) {
    $shp = $targetSlide.Shapes.Item($i)
    if ($idsToDelete -contains $shp.Id) {
        $shp.Delete()
    }
}
